$d = $word.ActiveDocument

# --- "Sample Courses" bullet: add a second line of Penn classes ---------
# The paragraph currently ends "... Crowdsourcing & Human Computation."
# Find that final phrase (including its trailing period) so we can turn
# the period into ", " and continue the sentence with a new italic run.
$rng = $d.Content
$found = $rng.Find.Execute("Crowdsourcing & Human Computation.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the match, then pull the start back one
    # character so the range covers only the trailing "."
    $rng.Collapse(0)
    $rng.MoveStart(1, -1) | Out-Null
    $rng.Text = ", "

    # Insert the second line of courses right after the ", " we just
    # wrote, as a new italic run (matching the style of the other course
    # names in this sentence).
    $ins = $rng.Duplicate
    $ins.Collapse(0)
    $insStart = $ins.Start
    $ins.InsertAfter("Technology & Policy, Theory of Networks, Intro to Dynamic Systems, Urban Education.") | Out-Null
    $ins.Italic = 1
    $ins.Font.NameBi = "Times New Roman"

    # Word keeps a "_GoBack" bookmark at the location of the most recent
    # edit. Re-adding it here (right between ", " and the new sentence)
    # both plants it at the correct spot and automatically removes the
    # stale one that used to sit over by "Appalachian Trail", since
    # bookmark names are unique and Bookmarks.Add replaces same-named
    # bookmarks.
    $bmRange = $d.Range($insStart, $insStart)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}
